# Daily attendance processing - 2025-10-02 15:04:17
# Updates recorded-by lists and attendance counts/percentages on the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (ANATOMY C1, session 2) ---
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("H3").Value = "108/221"

# --- Row 10 (Class Statistics - Average Attendance %) ---
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "33.1%"

# --- Row 12 (HISTOLOGY C1, session 1) ---
$ws.Range("G12").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 15 (Year 3 / C1 summary row) ---
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "38.9%"

# --- Row 16 (Year 3 / C2 summary row) ---
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = "27.2%"

# --- Row 25 (ANATOMY C2, session 2) ---
$ws.Range("G25").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("H25").Value = "63/246"

# --- Row 34 (HISTOLOGY C2, session 1) ---
$ws.Range("G34").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
